# "separate dept from affiliations"
#
# PI hours: split the old combined "dept" column (which actually held the
# list of affiliated units, e.g. "['ECE', 'CSL']") into two columns:
#   - dept: the PI's single home department (ECE / ME)
#   - app:  the full affiliation list (what "dept" used to hold)
#
# dept hours -> renamed to "department hours": now aggregated by the
# single home department only (ECE / ME), 2 rows.
#
# New sheet "unit(accumulative) hours" holds what the old "dept hours"
# sheet used to contain: hours/percentage accumulated across every unit
# each PI is affiliated with (CSL / ECE / ME / AE), 4 rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PI hours")
$ws2 = $wb.Worksheets.Item("dept hours")

# ---------------------------------------------------------------------
# 1. "PI hours": add an "app" column (F) holding the old affiliation
#    lists, and collapse the existing "dept" column (E) down to each
#    PI's single home department.
# ---------------------------------------------------------------------
$ws1.Range("F1").Value = "app"
$ws1.Range("E1").Copy()
$ws1.Range("F1").PasteSpecial(-4122)

$ws1.Range("F2").Value = $ws1.Range("E2").Value2
$ws1.Range("F3").Value = $ws1.Range("E3").Value2

$ws1.Range("E2").Value = "ECE"
$ws1.Range("E3").Value = "ME"

# ---------------------------------------------------------------------
# 2. Add the new "unit(accumulative) hours" sheet, right after
#    "dept hours", and seed it with the data currently in "dept hours"
#    (CSL/ECE/ME/AE accumulated hours) before that sheet is repurposed.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "unit(accumulative) hours"

$ws2.Range("B1:D1").Copy()
$ws3.Range("B1:D1").PasteSpecial(-4122)
$ws2.Range("A2:A5").Copy()
$ws3.Range("A2:A5").PasteSpecial(-4122)

$ws3.Range("B1").Value = "unit(accumulative)"
$ws3.Range("C1").Value = "hours"
$ws3.Range("D1").Value = "percentage"

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "CSL"
$ws3.Range("C2").Value = 10
$ws3.Range("D2").Value = 41.66666666666666

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "ECE"
$ws3.Range("C3").Value = 6
$ws3.Range("D3").Value = 25

$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "ME"
$ws3.Range("C4").Value = 4
$ws3.Range("D4").Value = 16.66666666666667

$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "AE"
$ws3.Range("C5").Value = 4
$ws3.Range("D5").Value = 16.66666666666667

# ---------------------------------------------------------------------
# 3. Rename "dept hours" to "department hours" and replace its data
#    with the per-home-department (ECE/ME) accumulation.
# ---------------------------------------------------------------------
$ws2.Name = "department hours"

$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(4).Delete()

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "ECE"
$ws2.Range("C2").Value = 6
$ws2.Range("D2").Value = 60

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "ME"
$ws2.Range("C3").Value = 4
$ws2.Range("D3").Value = 40

# Keep the original sheet ("PI hours") active/selected, matching the
# workbook's original tab selection.
$ws1.Activate()
